$d = $word.ActiveDocument

# --- Change 1: Methods paragraph -----------------------------------------
# "...regarding outlier removal and interpolation." ->
# "...regarding outlier removal, interpolation, and averaging methods."
$d.Content.Find.Execute(
    "were described in the methods section regarding outlier removal and interpolation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "were described in the methods section regarding outlier removal, interpolation, and averaging methods.",
    2)

# --- Change 2: Results paragraph ------------------------------------------
# Replace the sentence block (not touching the bold "Results:" heading run)
# with updated figures/wording.
$d.Content.Find.Execute(
    " Of the 7833 articles analyzed, 330 (4.2%) described their outlier removal procedures, and 472 (6.0%) described some interpolation procedures. The most popular outlier cutoffs are mean ± 3 or 4 SD (40.9% and 50.3%, respectively). When documented, the dominating interpolation time frame and procedure were one second (94.7%) and linear interpolation (92.8%), respectively.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Of the 7119 articles analyzed, 328 (4.6%) described outlier removal, 470 (6.6%) described interpolation, and an estimated 4366 (61.3%) described averaging methods. The most popular outlier cutoffs are mean ± 3 or 4 SD (40.9% and 50.3%, respectively). When documented, the dominating interpolation time frame and procedure were one second (94.7%) and linear interpolation (92.8%), respectively. Time-based bin averages (84.6%) were the most popular averaging methods.",
    2)
